$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")
Write-Host $ws.Range("A1").Value
Write-Host $ws.Name
Write-Host $wb.Worksheets.Count
